$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.04997025638596142
$ws.Range("D2").Value = 0.02157755715489174
$ws.Range("E2").Value = 0.424270739259299
$ws.Range("F2").Value = 0.382896176993512
$ws.Range("G2").Value = 0.2336590723292602
$ws.Range("H2").Value = 0.4056141175377235
$ws.Range("K2").Value = 1.78362691403629
$ws.Range("O2").Value = 1.193360450870998

$ws.Range("C3").Value = 0.04433884867539462
$ws.Range("D3").Value = 0.01901427093091712
$ws.Range("E3").Value = 0.3700728317306954
$ws.Range("F3").Value = 0.3838145461420837
$ws.Range("G3").Value = 0.2354554517334009
$ws.Range("H3").Value = 0.4115929439719608
$ws.Range("K3").Value = 1.557469793599068
$ws.Range("O3").Value = 1.209364920688103

$ws.Range("C4").Value = 0.04089501431946019
$ws.Range("D4").Value = 0.01743346039221905
$ws.Range("E4").Value = 0.3368859138142
$ws.Range("F4").Value = 0.3848471278943677
$ws.Range("G4").Value = 0.236951864078776
$ws.Range("H4").Value = 0.4156111646194205
$ws.Range("K4").Value = 1.418055232556583
$ws.Range("O4").Value = 1.220741970214036

$ws.Range("C5").Value = 0.03949509809710605
$ws.Range("D5").Value = 0.01678757195953295
$ws.Range("E5").Value = 0.3233829881485804
$ws.Range("F5").Value = 0.3853851691370451
$ws.Range("G5").Value = 0.2376598418865896
$ws.Range("H5").Value = 0.4173356763975136
$ws.Range("K5").Value = 1.36110720448039
$ws.Range("O5").Value = 1.225765901275011

$ws.Range("C6").Value = 0.03926285251057493
$ws.Range("D6").Value = 0.01668022178436246
$ws.Range("E6").Value = 0.3211420537095364
$ws.Range("F6").Value = 0.3854815747799591
$ws.Range("G6").Value = 0.2377833086799939
$ws.Range("H6").Value = 0.4176272807690751
$ws.Range("K6").Value = 1.351642954366298
$ws.Range("O6").Value = 1.226623474389186

$ws.Range("C7").Value = 0.04087612051199585
$ws.Range("D7").Value = 0.01742475651194297
$ws.Range("E7").Value = 0.336703725944858
$ws.Range("F7").Value = 0.3848539101436828
$ws.Range("G7").Value = 0.2369610155974442
$ws.Range("H7").Value = 0.4156340698455665
$ws.Range("K7").Value = 1.417287755151335
$ws.Range("O7").Value = 1.220808157791708

$ws.Range("C8").Value = 0.04802565737435316
$ws.Range("D8").Value = 0.02069520621004983
$ws.Range("E8").Value = 0.4055633357183694
$ws.Range("F8").Value = 0.3831151886882864
$ws.Range("G8").Value = 0.2341963547123385
$ws.Range("H8").Value = 0.4076034052631954
$ws.Range("K8").Value = 1.705764848277454
$ws.Range("O8").Value = 1.198555758758474

$ws.Range("C9").Value = 0.06215732992041012
$ws.Range("D9").Value = 0.02705157992572538
$ws.Range("E9").Value = 0.5414061261407994
$ws.Range("F9").Value = 0.3834521864023515
$ws.Range("G9").Value = 0.2319306385082669
$ws.Range("H9").Value = 0.3946214155329102
$ws.Range("K9").Value = 2.266947445502467
$ws.Range("O9").Value = 1.167315530999772

$ws.Range("C10").Value = 0.07261078708762625
$ws.Range("D10").Value = 0.0316848363551685
$ws.Range("E10").Value = 0.6418353489824966
$ws.Range("F10").Value = 0.3860229373528625
$ws.Range("G10").Value = 0.2322372491482554
$ws.Range("H10").Value = 0.3867853483643415
$ws.Range("K10").Value = 2.67635934335334
$ws.Range("O10").Value = 1.152054033099347

$ws.Range("C11").Value = 0.07738243928096722
$ws.Range("D11").Value = 0.03378424978525629
$ws.Range("E11").Value = 0.6876882950468826
$ws.Range("F11").Value = 0.3877053446988512
$ws.Range("G11").Value = 0.2328149684769869
$ws.Range("H11").Value = 0.3835933935306315
$ws.Range("K11").Value = 2.861959316267246
$ws.Range("O11").Value = 1.146809867433063

$ws.Range("C12").Value = 0.07919172017064113
$ws.Range("D12").Value = 0.03457801114149106
$ws.Range("E12").Value = 0.7050778580106822
$ws.Range("F12").Value = 0.3884168607559388
$ws.Range("G12").Value = 0.2330975767170713
$ws.Range("H12").Value = 0.3824385578004552
$ws.Range("K12").Value = 2.932145761991421
$ws.Range("O12").Value = 1.145070596478007

$ws.Range("C13").Value = 0.07880195400517209
$ws.Range("D13").Value = 0.03440711641856353
$ws.Range("E13").Value = 0.7013315138417653
$ws.Range("F13").Value = 0.3882603024356257
$ws.Range("G13").Value = 0.2330338596816546
$ws.Range("H13").Value = 0.3826848709700528
$ws.Range("K13").Value = 2.917034193402571
$ws.Range("O13").Value = 1.145434175282134

$ws.Range("C14").Value = 0.07753124254089983
$ws.Range("D14").Value = 0.03384957816689393
$ws.Range("E14").Value = 0.689118411256274
$ws.Range("F14").Value = 0.3877623861763126
$ws.Range("G14").Value = 0.2328369347108321
$ws.Range("H14").Value = 0.3834973028188529
$ws.Range("K14").Value = 2.867735547054451
$ws.Range("O14").Value = 1.146661821206663

$ws.Range("C15").Value = 0.07675320276933917
$ws.Range("D15").Value = 0.03350790679516535
$ws.Range("E15").Value = 0.6816409853450835
$ws.Range("F15").Value = 0.3874671100522065
$ws.Range("G15").Value = 0.2327246504752623
$ws.Range("H15").Value = 0.3840019673336883
$ws.Range("K15").Value = 2.83752606133362
$ws.Range("O15").Value = 1.147445971888885

$ws.Range("C16").Value = 0.07229927791647128
$ws.Range("D16").Value = 0.03154746383057727
$ws.Range("E16").Value = 0.6388422863420686
$ws.Range("F16").Value = 0.3859233667253577
$ws.Range("G16").Value = 0.2322083824424084
$ws.Range("H16").Value = 0.3870014750429789
$ws.Range("K16").Value = 2.664216681562607
$ws.Range("O16").Value = 1.152431128631576

$ws.Range("C17").Value = 0.06957112872288462
$ws.Range("D17").Value = 0.03034263908428869
$ws.Range("E17").Value = 0.6126309119016895
$ws.Range("F17").Value = 0.385108182544478
$ws.Range("G17").Value = 0.2320044960247074
$ws.Range("H17").Value = 0.3889372443523413
$ws.Range("K17").Value = 2.557729513797483
$ws.Range("O17").Value = 1.155926030153267

$ws.Range("C18").Value = 0.0680035050112906
$ws.Range("D18").Value = 0.02964887954424
$ws.Range("E18").Value = 0.5975704379714557
$ws.Range("F18").Value = 0.3846875528124158
$ws.Range("G18").Value = 0.231928426191736
$ws.Range("H18").Value = 0.3900857114642804
$ws.Range("K18").Value = 2.496420554108965
$ws.Range("O18").Value = 1.15809592224835

$ws.Range("C19").Value = 0.06747299790619365
$ws.Range("D19").Value = 0.02941385272422536
$ws.Range("E19").Value = 0.5924738414575756
$ws.Range("F19").Value = 0.384553400365867
$ws.Range("G19").Value = 0.2319097214561552
$ws.Range("H19").Value = 0.3904805760293399
$ws.Range("K19").Value = 2.475652140910427
$ws.Range("O19").Value = 1.158857970212708

$ws.Range("C20").Value = 0.06986138584741752
$ws.Range("D20").Value = 0.03047097544018129
$ws.Range("E20").Value = 0.6154195268870097
$ws.Range("F20").Value = 0.3851899625615687
$ws.Range("G20").Value = 0.2320219300708359
$ws.Range("H20").Value = 0.3887275473102676
$ws.Range("K20").Value = 2.569071528298991
$ws.Range("O20").Value = 1.15553744470887

$ws.Range("C21").Value = 0.07790441715584961
$ws.Range("D21").Value = 0.03401337474402055
$ws.Range("E21").Value = 0.6927049725179444
$ws.Range("F21").Value = 0.3879066111061462
$ws.Range("G21").Value = 0.2328930372540441
$ws.Range("H21").Value = 0.3832572070406997
$ws.Range("K21").Value = 2.882218395459631
$ws.Range("O21").Value = 1.146294520906963

$ws.Range("C22").Value = 0.08317476358521958
$ws.Range("D22").Value = 0.03632128143586044
$ws.Range("E22").Value = 0.7433682505728001
$ws.Range("F22").Value = 0.3901162103968119
$ws.Range("G22").Value = 0.2338348364110345
$ws.Range("H22").Value = 0.379996262071927
$ws.Range("K22").Value = 3.086315367057296
$ws.Range("O22").Value = 1.141692093130501

$ws.Range("C23").Value = 0.08036061718823362
$ws.Range("D23").Value = 0.03509018860589208
$ws.Range("E23").Value = 0.7163136582145029
$ws.Range("F23").Value = 0.3888969656880192
$ws.Range("G23").Value = 0.2332978274802144
$ws.Range("H23").Value = 0.3817078399737994
$ws.Range("K23").Value = 2.977437655808558
$ws.Range("O23").Value = 1.144016106068989

$ws.Range("C24").Value = 0.06973015803289684
$ws.Range("D24").Value = 0.03041295796452914
$ws.Range("E24").Value = 0.6141587668514745
$ws.Range("F24").Value = 0.3851528402507185
$ws.Range("G24").Value = 0.2320139200274625
$ws.Range("H24").Value = 0.3888222405689419
$ws.Range("K24").Value = 2.563944077599274
$ws.Range("O24").Value = 1.155712624008856

$ws.Range("C25").Value = 0.0583220997053786
$ws.Range("D25").Value = 0.02533832570496486
$ws.Range("E25").Value = 0.5045568252794084
$ws.Range("F25").Value = 0.382955873145697
$ws.Range("G25").Value = 0.2322008294922568
$ws.Range("H25").Value = 0.3978356557814706
$ws.Range("K25").Value = 2.115630613574581
$ws.Range("O25").Value = 1.145434175282134

Write-Output "Updated 192 cells"